# Applies the commit "sampai labeling manual, sesuai alur baru":
# appends 15 new crawled-tweet rows (111-125) to Sheet1, updates the
# used dimension/sheet view accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(111, 1).Value = 1339438533426044928
$ws.Cells.Item(111, 2).Value = 'Ini loh Mugwort yang lagi jadi holy grail-nya sobat tiktok sekalian karena emang hasilnya sebagus itu di muka :'')
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #mugwort https://t.co/kArzlIgUUJ'
$ws.Cells.Item(111, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(111, 4).Value = 'Thu Dec 17 05:13:26 +0000 2020'

$ws.Cells.Item(112, 1).Value = 1339384794988257024
$ws.Cells.Item(112, 2).Value = 'Nilai aja yang tinggi nggak ada yang masuk di otak 
#belajaronline'
$ws.Cells.Item(112, 3).Value = 'AbdSalam24'
$ws.Cells.Item(112, 4).Value = 'Thu Dec 17 01:39:54 +0000 2020'

$ws.Cells.Item(113, 1).Value = 1339082840617598976
$ws.Cells.Item(113, 2).Value = 'Yuk join ke Groupnya untuk informasi setiap minggunya di https://t.co/xjVIJXwfwr
Bonjour Prancis
One Third Consulting &amp; Abroad
#BonjourPrancis #lesbahasaprancis #konsultasipendidikan #belajaronline #beasiswa #beasiswaprancis #konsultasigratis #beasiswaeropa #bahasaprancis https://t.co/KZPrhCbrAs'
$ws.Cells.Item(113, 3).Value = 'BonjourPrancis'
$ws.Cells.Item(113, 4).Value = 'Wed Dec 16 05:40:02 +0000 2020'

$ws.Cells.Item(114, 1).Value = 1339066669096825088
$ws.Cells.Item(114, 2).Value = 'Ingat lagu ini? 
.
Sc twitter @tvindonesiawkwk
.
#kampuscenter #workshoplulusptn #PJJ #belajaronline #sma #kampus #pelajar #ltmpt #kuliah #UTBK #SBMPTN #snmptn2020  #ltmpt2020 #ltmpt #sbmptn2020 #utbk2020 #utbk #masukkampus #kampusindonesia https://t.co/v8g762UPK5'
$ws.Cells.Item(114, 3).Value = 'kampuscenter'
$ws.Cells.Item(114, 4).Value = 'Wed Dec 16 04:35:47 +0000 2020'

$ws.Cells.Item(115, 1).Value = 1338861939586298112
$ws.Cells.Item(115, 2).Value = 'Hari ini di kelas Jena kita membahas soal ujian minggu kemaren😁 Udah pada jago nih temen-temen, jangan lupa kudu PD ya buat sprechen nyaa👩‍💻👨‍💻
#kursusbahasa #kursusonline #belajarjerman #belajaronline #ujian https://t.co/Wla9l8laOD'
$ws.Cells.Item(115, 3).Value = 'jermanstudiere'
$ws.Cells.Item(115, 4).Value = 'Tue Dec 15 15:02:15 +0000 2020'

$ws.Cells.Item(116, 1).Value = 1338750899888374016
$ws.Cells.Item(116, 2).Value = 'Masih susah jamanku ya...
.
Sc: unknown (yang tau kasih tau bro)
.
.
Reposted from @1cak
.
#workshoplulusptn #kampus #kampuscenter #belajaronline #sma #kampus #pelajar #ltmpt #kuliah #UTBK #SNMPTN #SBMPTN #snmptn2020  #ltmpt2020 #ltmpt #sbmptn2020 #utbk2020 #utbk #masukkampus https://t.co/HtJ87PPW5j'
$ws.Cells.Item(116, 3).Value = 'kampuscenter'
$ws.Cells.Item(116, 4).Value = 'Tue Dec 15 07:41:01 +0000 2020'

$ws.Cells.Item(117, 1).Value = 1338721025832361984
$ws.Cells.Item(117, 2).Value = 'Terima kasih, suhu~ @MelissaSunjaya 
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisa https://t.co/gTFH9khSL7'
$ws.Cells.Item(117, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(117, 4).Value = 'Tue Dec 15 05:42:19 +0000 2020'

$ws.Cells.Item(118, 1).Value = 1338720825596236032
$ws.Cells.Item(118, 2).Value = 'Ya, gimana ya. Aduh mba Mel, aku bingung menjelaskannya~
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisa https://t.co/jYQwJNvBJr'
$ws.Cells.Item(118, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(118, 4).Value = 'Tue Dec 15 05:41:31 +0000 2020'

$ws.Cells.Item(119, 1).Value = 1338720603189125120
$ws.Cells.Item(119, 2).Value = 'Iya plissss pulsa internet sama subscribe layanan streaming-nya sekalian~
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisa https://t.co/d0b1vfWkvc'
$ws.Cells.Item(119, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(119, 4).Value = 'Tue Dec 15 05:40:38 +0000 2020'

$ws.Cells.Item(120, 1).Value = 1338720246790716928
$ws.Cells.Item(120, 2).Value = 'Ga lengkap kalo ga lucu-lucuan, ya kan?!
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisa https://t.co/vI793G1wc6'
$ws.Cells.Item(120, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(120, 4).Value = 'Tue Dec 15 05:39:13 +0000 2020'

$ws.Cells.Item(121, 1).Value = 1338720054427336960
$ws.Cells.Item(121, 2).Value = 'Pas banget buat sobi twitter yang hobinya 2 a.m thoughts :")
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisa https://t.co/nH1IVEfeHk'
$ws.Cells.Item(121, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(121, 4).Value = 'Tue Dec 15 05:38:27 +0000 2020'

$ws.Cells.Item(122, 1).Value = 1338719595423621120
$ws.Cells.Item(122, 2).Value = 'Yang udah dapet jadwal tuker kado akhir tahun, sini dulu yuk!
Ada rekomendasi kado paling ''sumpah ga ngerti lagi'' dari Melissa Sunjaya biar momen tuker kado kamu lebih asik, lop.
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #tulisan https://t.co/od6zPaYitb'
$ws.Cells.Item(122, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(122, 4).Value = 'Tue Dec 15 05:36:38 +0000 2020'

$ws.Cells.Item(123, 1).Value = 1338692370167042048
$ws.Cells.Item(123, 2).Value = 'Mau tahu cara praktis dan murah untuk Pembelajaran Jarak Jauh (PJJ)?👩‍💻
Pasti mau kan... Saatnya kamu tinggalkan cara manual untuk mengatur jadwal kelas dan meeting secara virtual 📆.
Solusinya, pakai https://t.co/1rIQ8ab0xT.
#belajaronline 
#fingerspot https://t.co/SlWG2pIvm4'
$ws.Cells.Item(123, 3).Value = 'Bambang90500845'
$ws.Cells.Item(123, 4).Value = 'Tue Dec 15 03:48:27 +0000 2020'

$ws.Cells.Item(124, 1).Value = 1338404788174606080
$ws.Cells.Item(124, 2).Value = '- di twitter gw belajar pacaran online tp sakitnya beneran
- di twitter gw belajar jembud cewe dipanggang
- di twitter gw belajar ngetik wkwkwk tanpa ketawa
- di twitter gw belajar foto dikaca toilet mall

#belajaronline'
$ws.Cells.Item(124, 3).Value = 'josephine_oi'
$ws.Cells.Item(124, 4).Value = 'Mon Dec 14 08:45:42 +0000 2020'

$ws.Cells.Item(125, 1).Value = 1338391443543532032
$ws.Cells.Item(125, 2).Value = 'Karena punya ide tanpa aksi sama dengan kosong. Seperti kata tukang sepatu terkenal, "just do it!". Ingat, perjuangan tidak akan pernah mengkhianati hasil.

Selamat hari Senin!

#kelaskita #carabarubelajarseru #belajardirumah
#elearning #belajaronline #dirumahaja https://t.co/BIotzSic8C'
$ws.Cells.Item(125, 3).Value = 'kelaskitadotcom'
$ws.Cells.Item(125, 4).Value = 'Mon Dec 14 07:52:40 +0000 2020'

$ws.Application.ActiveWindow.ScrollRow = 103
$ws.Range("J121").Select()

Write-Output "applied rows 111-125"
